# Publication 0.2.0 preparation
# - Bump Version from 0.1.1 to 0.2.0
# - Bump Date to the new publication timestamp
# - Insert a new "Jurisdiction" metadata row (after "Contact")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version (row 3) and Date (row 8) values
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# Insert a new row right after "Contact" (row 10), pushing Description.. down by one
$ws.Rows.Item(11).Insert()

# Carry over the same formatting used by the other metadata rows
$ws.Range("A10:B10").Copy($ws.Range("A11:B11"))

# Populate the new Jurisdiction row
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
